$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1239.0714
$ws.Range("I17").Value = 666
$ws.Range("J17").Value = 1363.6522
$ws.Range("K17").Value = 1998
$ws.Range("L17").Value = 4090.9566
$ws.Range("M17").Value = -1830
$ws.Range("N17").Value = -4426.9566

$ws.Range("H28").Value = 1248.75
$ws.Range("I28").Value = 550.9286
$ws.Range("J28").Value = 2877
$ws.Range("K28").Value = 550.9286
$ws.Range("L28").Value = 2877
$ws.Range("M28").Value = -65.92859999999996
$ws.Range("N28").Value = -3847

$ws.Range("H53").Value = 280.9375
$ws.Range("I53").Value = 276
$ws.Range("J53").Value = 283.9
$ws.Range("K53").Value = 276
$ws.Range("L53").Value = 283.9
$ws.Range("M53").Value = 361
$ws.Range("N53").Value = -1557.9

$ws.Range("H87").Value = 49998
$ws.Range("J87").Value = 49998
$ws.Range("L87").Value = 49998
$ws.Range("N87").Value = -52494

$ws.Range("H90").Value = 49998
$ws.Range("J90").Value = 49998
$ws.Range("L90").Value = 149994
$ws.Range("N90").Value = -162474

$ws.Range("H135").Value = 2001390
$ws.Range("I135").Value = 3334983.2
$ws.Range("K135").Value = 30014848.8
$ws.Range("M135").Value = -30012313.8

$ws.Range("H137").Value = 4349.0527
$ws.Range("I137").Value = 3116.8635
$ws.Range("J137").Value = 6043.3125
$ws.Range("K137").Value = 9350.5905
$ws.Range("L137").Value = 18129.9375
$ws.Range("M137").Value = -6800.5905
$ws.Range("N137").Value = -23229.9375

$ws.Range("H138").Value = 2226888.2
$ws.Range("J138").Value = 4173840
$ws.Range("L138").Value = 12521520
$ws.Range("N138").Value = -12531800

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H21").Value = 623
$ws.Range("I21").Value = 623
$ws.Range("K21").Value = 623
$ws.Range("M21").Value = -249

$ws.Range("H30").Value = 500
$ws.Range("I30").Value = 500
$ws.Range("K30").Value = 500
$ws.Range("M30").Value = -350

$ws.Range("H61").Value = 8625.764999999999
$ws.Range("I61").Value = 1327.875
$ws.Range("K61").Value = 1327.875
$ws.Range("M61").Value = -1115.875

$ws.Range("H97").Value = 2381778.5
$ws.Range("I97").Value = 612.86206
$ws.Range("J97").Value = 13890746
$ws.Range("K97").Value = 612.86206
$ws.Range("L97").Value = 13890746
$ws.Range("M97").Value = -116.86206
$ws.Range("N97").Value = -13891738

$ws.Range("H132").Value = 5803.8223
$ws.Range("I132").Value = 3555.182
$ws.Range("K132").Value = 10665.546
$ws.Range("M132").Value = -8135.545999999998

$ws.Range("H136").Value = 8625.764999999999
$ws.Range("I136").Value = 1327.875
$ws.Range("K136").Value = 3983.625
$ws.Range("M136").Value = -1433.625

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 4506167.5
$ws.Range("I20").Value = 6411573
$ws.Range("K20").Value = 6411573
$ws.Range("M20").Value = -6411326

$ws.Range("H94").Value = 2627.2144
$ws.Range("I94").Value = 1814.6666
$ws.Range("K94").Value = 1814.6666
$ws.Range("M94").Value = -1363.6666

$ws.Range("H99").Value = 5349586
$ws.Range("I99").Value = 2113.111
$ws.Range("J99").Value = 11365494
$ws.Range("K99").Value = 2113.111
$ws.Range("L99").Value = 11365494
$ws.Range("M99").Value = -615.1109999999999
$ws.Range("N99").Value = -11368490

$ws.Range("H132").Value = 109333.336
$ws.Range("J132").Value = 109333.336
$ws.Range("L132").Value = 109333.336
$ws.Range("N132").Value = -119453.336

$ws.Range("H134").Value = 7719.6294
$ws.Range("J134").Value = 10363.111
$ws.Range("L134").Value = 31089.333
$ws.Range("N134").Value = -36159.333

$ws.Range("H135").Value = 79994.5
$ws.Range("J135").Value = 79994.5
$ws.Range("L135").Value = 79994.5
$ws.Range("N135").Value = -90134.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 66666740
$ws.Range("I7").Value = 40.8
$ws.Range("J7").Value = 200000140
$ws.Range("K7").Value = 40.8
$ws.Range("L7").Value = 200000140
$ws.Range("M7").Value = 72.2
$ws.Range("N7").Value = -200000366

$ws.Range("H31").Value = 7747.6665
$ws.Range("I31").Value = 2551.9412
$ws.Range("K31").Value = 2551.9412
$ws.Range("M31").Value = -2256.9412

$ws.Range("H34").Value = 7747.6665
$ws.Range("I34").Value = 2551.9412
$ws.Range("K34").Value = 2551.9412
$ws.Range("M34").Value = -2349.9412

$ws.Range("H86").Value = 113669220
$ws.Range("J86").Value = 333333340
$ws.Range("L86").Value = 333333340
$ws.Range("N86").Value = -333335586

$ws.Range("H89").Value = 113669220
$ws.Range("J89").Value = 333333340
$ws.Range("L89").Value = 1666666700
$ws.Range("N89").Value = -1666677932

$ws.Range("H122").Value = 4300.4814
$ws.Range("I122").Value = 2921.7693
$ws.Range("K122").Value = 8765.3079
$ws.Range("M122").Value = -6315.3079

$ws.Range("H132").Value = 4523.9067
$ws.Range("I132").Value = 2357.3333
$ws.Range("J132").Value = 7260.6313
$ws.Range("K132").Value = 7071.999899999999
$ws.Range("L132").Value = 21781.8939
$ws.Range("M132").Value = -4541.999899999999
$ws.Range("N132").Value = -26841.8939

$ws.Range("H134").Value = 7465.8667
$ws.Range("I134").Value = 3430.8333
$ws.Range("K134").Value = 10292.4999
$ws.Range("M134").Value = -7757.499899999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 33337294
$ws.Range("I14").Value = 33337294
$ws.Range("K14").Value = 100011882
$ws.Range("M14").Value = -100011709

$ws.Range("I38").Value = 37.333332
$ws.Range("K38").Value = 111.999996
$ws.Range("M38").Value = 235.000004

$ws.Range("H69").Value = 5937.6665

$ws.Range("H72").Value = 5937.6665

$ws.Range("H100").Value = 5360
$ws.Range("J100").Value = 5277.5
$ws.Range("L100").Value = 15832.5
$ws.Range("N100").Value = -17454.5

$ws.Range("H102").Value = 15000
$ws.Range("J102").Value = 15000
$ws.Range("L102").Value = 45000
$ws.Range("N102").Value = -49868

$ws.Range("H104").Value = 6983.1665
$ws.Range("J104").Value = 7499.75
$ws.Range("L104").Value = 22499.25
$ws.Range("N104").Value = -27741.25

$ws.Range("H105").Value = 14995
$ws.Range("J105").Value = 14995
$ws.Range("L105").Value = 44985
$ws.Range("N105").Value = -50227

$ws.Range("H107").Value = 18183058
$ws.Range("I107").Value = 686.3333
$ws.Range("J107").Value = 25001446
$ws.Range("K107").Value = 2058.9999
$ws.Range("L107").Value = 75004338
$ws.Range("M107").Value = -138.9998999999998
$ws.Range("N107").Value = -75008178

$ws.Range("H113").Value = 5538.6
$ws.Range("J113").Value = 6042.8887
$ws.Range("L113").Value = 18128.6661
$ws.Range("N113").Value = -22468.6661

$ws.Range("H139").Value = 95961.91
$ws.Range("I139").Value = 128822.625
$ws.Range("K139").Value = 386467.875
$ws.Range("M139").Value = -381327.875

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H22").Value = 369.6
$ws.Range("I22").Value = 369.6
$ws.Range("K22").Value = 369.6
$ws.Range("M22").Value = 159.4

$ws.Range("H46").Value = 40811.11
$ws.Range("J46").Value = 56383.332
$ws.Range("L46").Value = 56383.332
$ws.Range("N46").Value = -56695.332

$ws.Range("H102").Value = 3065.8064
$ws.Range("I102").Value = 3107.8845
$ws.Range("J102").Value = 2847
$ws.Range("K102").Value = 3107.8845
$ws.Range("L102").Value = 2847
$ws.Range("M102").Value = -1485.8845
$ws.Range("N102").Value = -6091

$ws.Range("H113").Value = 309675.5
$ws.Range("I113").Value = 836049.5
$ws.Range("K113").Value = 836049.5
$ws.Range("M113").Value = -833879.5

$ws.Range("H126").Value = 33339988
$ws.Range("I126").Value = 100002400
$ws.Range("K126").Value = 300007200
$ws.Range("M126").Value = -300004730

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H56").Value = 34500
$ws.Range("I56").Value = 34500
$ws.Range("K56").Value = 34500
$ws.Range("M56").Value = -33809

$ws.Range("H107").Value = 2899.3333
$ws.Range("I107").Value = 2899.3333
$ws.Range("K107").Value = 2899.3333
$ws.Range("M107").Value = -979.3332999999998

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H61").Value = 2528
$ws.Range("J61").Value = 2057
$ws.Range("L61").Value = 2057
$ws.Range("N61").Value = -2641

$ws.Range("H96").Value = 865.5
$ws.Range("I96").Value = 898.6
$ws.Range("K96").Value = 898.6
$ws.Range("M96").Value = 474.4

$ws.Range("H132").Value = 43485470
$ws.Range("I132").Value = 66674388
$ws.Range("K132").Value = 200023164
$ws.Range("M132").Value = -200020634

$ws.Range("H136").Value = 404863.6
$ws.Range("I136").Value = 1152.4546
$ws.Range("K136").Value = 3457.3638
$ws.Range("M136").Value = -907.3638000000001

